try {
    $Error.Clear()
    $bytes = [System.IO.File]::ReadAllBytes("C:\nonexistent.pptx")
    Write-Output "read ok: $($bytes.Length)"
} catch {
    Write-Output "error: $_"
}
